$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: rename "LEGAL FEES" -> "LEGAL / BAILIFF FEES" ---
# (done first so its new shared-string slot lands right after FRANCE,
#  matching the order the fixture data was originally authored in)
$ws.Range("A8").Value = "LEGAL / BAILIFF FEES"

# --- Header row (row 1): rename columns ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Parent"
$ws.Range("C1").Value = "AtPath"

# --- Rows 3-13: add the AtPath value for every charge row ---
# Leading apostrophe reproduces the "quote prefix" cell formatting Excel
# applies to text values that could otherwise be misread (e.g. "/FRA").
for ($r = 3; $r -le 13; $r++) {
    $ws.Range("C$r").Value = "'/FRA"
}

# --- Row 2 (FRANCE section row): add the AtPath value ---
$ws.Range("C2").Value = "'/FRA"

# --- Update the active selection from B4 to C4 ---
$ws.Range("C4").Select()
